$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- 1. Update the URL value (pythia -> cicada) ---
$ws.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/VaccineGender"

# --- 2. Update the Date value ---
$ws.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# --- 3. Insert a new "Jurisdiction" row before "Description" (currently row 11) ---
# Shift rows 11..21 down into 12..22 (bottom-up to avoid clobbering source data)
for ($r = 21; $r -ge 11; $r--) {
    $destRow = $r + 1
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value2
}

# Make sure the new last row (22) carries the same formatting as the row above it
$ws.Range("A21:B21").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)

# Populate the newly freed row 11 with the Jurisdiction property
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
